$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '90.372.26'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -0.27%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.092.87'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -1.87%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '233.75'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +8.82%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '624.50'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.35%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -6.08%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.363'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.83%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +0.07%  '
$ws.Range('B10').NumberFormat = '@'
$ws.Range('B10').Value = 'LidoStakedEther'
$ws.Range('C10').NumberFormat = '@'
$ws.Range('C10').Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '3.091.54'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -1.81%  '
$ws.Range('B11').NumberFormat = '@'
$ws.Range('B11').Value = 'Cardano'
$ws.Range('C11').NumberFormat = '@'
$ws.Range('C11').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.729'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -8.45%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.197'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -1.69%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '36.63'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000253'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +3.31%  '
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -3.98%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '90.033.10'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -0.48%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.652.55'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -2.25%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.080.73'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -3.05%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.80'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +3.51%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0000214'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +1.05%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.98'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -2.09%  '
$ws.Range('B22').NumberFormat = '@'
$ws.Range('B22').Value = 'BitcoinCash'
$ws.Range('C22').NumberFormat = '@'
$ws.Range('C22').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '436.69'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -5.10%  '
$ws.Range('B23').NumberFormat = '@'
$ws.Range('B23').Value = 'Polkadot'
$ws.Range('C23').NumberFormat = '@'
$ws.Range('C23').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.55'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +6.49%  '
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -0.89%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '5.91'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -0.13%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -1.48%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '88.59'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -3.70%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '12.22'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +0.67%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '3.250.58'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -2.20%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +0.00%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '9.42'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +2.40%  '
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -1.58%  '
$ws.Range('B33').NumberFormat = '@'
$ws.Range('B33').Value = 'Stellar'
$ws.Range('C33').NumberFormat = '@'
$ws.Range('C33').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.199'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +8.07%  '
$ws.Range('B34').NumberFormat = '@'
$ws.Range('B34').Value = 'EthereumClassic'
$ws.Range('C34').NumberFormat = '@'
$ws.Range('C34').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '25.92'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -2.59%  '
$ws.Range('B35').NumberFormat = '@'
$ws.Range('B35').Value = 'Kaspa'
$ws.Range('C35').NumberFormat = '@'
$ws.Range('C35').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.154'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +6.57%  '
$ws.Range('B36').NumberFormat = '@'
$ws.Range('B36').Value = 'dogwifhat'
$ws.Range('C36').NumberFormat = '@'
$ws.Range('C36').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.87'
$ws.Range('B37').NumberFormat = '@'
$ws.Range('B37').Value = 'Binance-PegBSC-USD'
$ws.Range('C37').NumberFormat = '@'
$ws.Range('C37').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.894'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -10.82%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '505.34'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -2.32%  '
$ws.Range('B39').NumberFormat = '@'
$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').NumberFormat = '@'
$ws.Range('C39').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '7.03'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +2.06%  '
$ws.Range('B40').NumberFormat = '@'
$ws.Range('B40').Value = 'PancakeSwap'
$ws.Range('C40').NumberFormat = '@'
$ws.Range('C40').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.90'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -0.41%  '
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -0.20%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0870'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -1.34%  '
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -0.12%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.409'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -2.67%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.49'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +55.40%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.90'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -1.68%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '151.24'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +2.56%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.690'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +3.77%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '44.90'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +0.87%  '
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +0.75%  '
